# Weekly update: insert a new data row at row 20 (pushing the existing
# rows 20-30 down to 21-31) and populate the new row with this week's
# Ciboulette price record for Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20:30 down to 21:31, leaving row 20 free for the new record.
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with the new weekly record.
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 45215
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = 100112039
$ws.Cells.Item(20, 7).Value = "Ciboulette"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = 2000
$ws.Cells.Item(20, 14).Value = "`$/docena de atados"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 667
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = "Hortaliza"
